# KnightTraverse function in Knightmove.h DONE with backtracking
# Applies the TimeScheme.xlsx edits: new task-session rows, a few corrected
# timestamps, three new blank template rows (12-14) before the totals row,
# and the totals/footer rows shifting down accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert three blank rows above the old "totals" row (row 12), pushing
#    the totals row and the two footer rows down to 15/16/17. Excel keeps
#    all formulas (SUM ranges, $D$14 absolute refs, etc.) consistent.
# ---------------------------------------------------------------------
$ws.Rows("12:14").Insert()

# ---------------------------------------------------------------------
# 2) Corrected Start/Stop timestamps
# ---------------------------------------------------------------------
$ws.Range("C2").Value = 41294.680844907409
$ws.Range("B11").Value = 41294.680844907409

# ---------------------------------------------------------------------
# 3) Fill in Start/Stop timestamps for rows 5, 6 and 7 (previously blank)
# ---------------------------------------------------------------------
$ws.Range("B4").Copy()
$ws.Range("B5:C7").PasteSpecial(-4122)

$ws.Range("B5").Value = 41294.736635879628
$ws.Range("C5").Value = 41294.740485416667
$ws.Range("B6").Value = 41294.740485416667
$ws.Range("C6").Value = 41294.744398148148
$ws.Range("B7").Value = 41294.815879629627
$ws.Range("C7").Value = 41294.947791319442

# ---------------------------------------------------------------------
# 4) Row 11: remove the task-budget (E11) and the "Comments" (J11) values
# ---------------------------------------------------------------------
$ws.Range("E11").ClearContents()
$ws.Range("J11").ClearContents()

# ---------------------------------------------------------------------
# 5) "Comments" column: mark rows 2,3,4,6 as "Good time"
# ---------------------------------------------------------------------
$ws.Range("J2").Value = "Good time"
$ws.Range("J3").Value = "Good time"
$ws.Range("J4").Value = "Good time"
$ws.Range("J6").Value = "Good time"

# ---------------------------------------------------------------------
# 6) Fill formulas for the new template rows 12-14 (columns D, F, G, H, I)
# ---------------------------------------------------------------------
$ws.Range("D12:D14").Formula = '=HOUR(C12)+MINUTE(C12)/60-(HOUR(B12)+MINUTE(B12)/60)'
$ws.Range("F12:F14").Value = 0
$ws.Range("G12:G14").Formula = '=D12*$D$17'
$ws.Range("H12:H14").Formula = '=IF(G12>0,F12,0)'
$ws.Range("I12:I14").Formula = '=IF(G12>0,H12+H12-G12,0)'
$ws.Range("G12:I14").Style = "Normal"

# ---------------------------------------------------------------------
# 7) Totals row (now row 15) needs to sum through the new rows (H14/I14)
# ---------------------------------------------------------------------
$ws.Range("H15").Formula = '=SUM(H2:H14)'
$ws.Range("I15").Formula = '=SUM(I2:I14)'

# ---------------------------------------------------------------------
# 8) Selection marker, as recorded for this edit
# ---------------------------------------------------------------------
$ws.Range("I16").Select()
